# Fill in seller/buyer invoice lookup details for rows 2, 4, 5, 7, 8, 10
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns E:L are treated as text so numeric-looking values
# (e.g. invoice numbers / tax codes with leading zeros) are preserved.
# Only the rows that actually receive new data are touched.
$ws.Range("E2:L2").NumberFormat = "@"
$ws.Range("E4:L4").NumberFormat = "@"
$ws.Range("E5:L5").NumberFormat = "@"
$ws.Range("E7:L7").NumberFormat = "@"
$ws.Range("E8:L8").NumberFormat = "@"
$ws.Range("E10:L10").NumberFormat = "@"

# Row 2
$ws.Range("E2").Value = "00068496"
$ws.Range("F2").Value = "CÔNG TY CỔ PHẦN BA HUÂN"
$ws.Range("G2").Value = "0304244470"
$ws.Range("H2").Value = "22 Nguyễn Đình Chi, Phường 09, Quận 6, Thành phố Hồ Chí Minh, Việt Nam"
$ws.Range("J2").Value = "Liên Hiệp Hợp Tác Xã Thương Mại TP. Hồ Chí Minh"
$ws.Range("K2").Value = "199-205 Nguyễn Thái Học, Phường Phạm Ngũ Lão, Quận 1, Thành phố Hồ Chí Minh, Việt Nam"
$ws.Range("L2").Value = "0301175691"

# Row 4
$ws.Range("E4").Value = "00023079"
$ws.Range("F4").Value = "CÔNG TY CỔ PHẦN PHÚ TRƯỜNG QUỐC TẾ"
$ws.Range("G4").Value = "0304308445"
$ws.Range("H4").Value = "15A1 Đường Nguyễn Hữu Thọ, Xã Phước Kiển, Huyện Nhà Bè, Thành phố Hồ Chí Minh, Việt Nam"
$ws.Range("I4").Value = "0181003527080"
$ws.Range("J4").Value = "Liên Hiệp Hợp Tác Xã Thương Mại TP.HCM"
$ws.Range("K4").Value = "199-205 Nguyễn Thái Học, Phường Phạm Ngũ Lão, Quận 01, Thành Phố Hồ Chí Minh"
$ws.Range("L4").Value = "0301175691"

# Row 5
$ws.Range("E5").Value = "00002426"
$ws.Range("F5").Value = "CÔNG TY TNHH CHĂN NUÔI TAFA VIỆT"
$ws.Range("G5").Value = "3401142134"
$ws.Range("H5").Value = "Thôn 1, Xã Trà Tân, Huyện Đức Linh, Tỉnh Bình Thuận, Việt Nam"
$ws.Range("I5").Value = "5400 201 010 928"
$ws.Range("J5").Value = "CÔNG TY TNHH MỘT THÀNH VIÊN TMDV SIÊU THỊ CO.OPMART ĐÀ NẴNG"
$ws.Range("K5").Value = "478 Điện Biên Phủ, Phường Thanh Khê Đông, Quận Thanh Khê, Thành phố Đà Nẵng, Việt Nam"
$ws.Range("L5").Value = "0401281414"

# Row 7
$ws.Range("E7").Value = "00000094"
$ws.Range("F7").Value = "CÔNG TY TNHH ĐẦU TƯ - SẢN XUẤT VÀ THƯƠNG MẠI HOÀNG BÁCH"
$ws.Range("G7").Value = "0318580994"
$ws.Range("H7").Value = "C9/16A31 Bùi Thanh Khiết, khu phố 3, Thị Trấn Tân Túc, Huyện Bình Chánh, Thành phố Hồ Chí Minh, Việt Nam"
$ws.Range("I7").Value = "110605110688"
$ws.Range("J7").Value = "CHI NHÁNH LIÊN HIỆP HỢP TÁC XÃ THƯƠNG MẠI TP.HỒ CHÍ MINH - CO.OPMART CHU VĂN AN"
$ws.Range("K7").Value = "Tầng 1-Tầng 2 Khối A&B Cao ốc Đất Phương Nam 241A Chu Văn An, Phường 12, Quận Bình Thạnh, TP Hồ Chí Minh"
$ws.Range("L7").Value = "0301175691-036"

# Row 8
$ws.Range("E8").Value = "00004221"
$ws.Range("F8").Value = "CÔNG TY TNHH SẢN XUẤT THƯƠNG MẠI DỊCH VỤ NAM KHẢI PHÚ"
$ws.Range("G8").Value = "0312563329"
$ws.Range("H8").Value = "81 Cách Mạng Tháng Tám, Phường Bến Thành, Quận 1, Hồ Chí Minh"
$ws.Range("J8").Value = "Chi Nhánh Liên Hiệp Hợp Tác Xã Thương Mại TP.Hồ Chí Minh - Co.opMart Cái Bè"
$ws.Range("K8").Value = "Khu 2, Thị trấn Cái Bè, Huyện Cái Bè, Tỉnh Tiền Giang"
$ws.Range("L8").Value = "0301175691-068"

# Row 10
$ws.Range("E10").Value = "00004288"
$ws.Range("F10").Value = "CÔNG TY TNHH SẢN XUẤT THƯƠNG MẠI DỊCH VỤ NAM KHẢI PHÚ"
$ws.Range("G10").Value = "0312563329"
$ws.Range("H10").Value = "81 Cách Mạng Tháng Tám, Phường Bến Thành, Quận 1, Hồ Chí Minh"
$ws.Range("J10").Value = "Chi Nhánh Liên Hiệp Hợp Tác Xã Thương Mại TP.Hồ Chí Minh - Co.opMart Cái Bè"
$ws.Range("K10").Value = "Khu 2, Thị trấn Cái Bè, Huyện Cái Bè, Tỉnh Tiền Giang"
$ws.Range("L10").Value = "0301175691-068"
